$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Window width change
$excel.ActiveWindow.Width = 27945

# Update frozen pane top-left cell and active selection
$ws.Range("E13").Select()
$ws.Range("N34").Select()

# D31: replace shared formula with explicit formula (value unchanged: 243)
$ws.Range("D31").Formula = "=B31-C31"

# Row 32 data
$ws.Range("B32").Value = 1652
$ws.Range("C32").Value = 1369
$ws.Range("D32").Formula = "=B32-C32"
$ws.Range("E32").Value = 61
$ws.Range("F32").Value = 6
$ws.Range("G32").Value = 7
$ws.Range("H32").Value = 334.26
$ws.Range("I32").Value = 3
$ws.Range("J32").Value = 329

# Row 33 data
$ws.Range("B33").Value = 995
$ws.Range("C33").Value = 799
$ws.Range("D33").Formula = "=B33-C33"
$ws.Range("E33").Value = 33
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 219.13
$ws.Range("I33").Value = 1
$ws.Range("J33").Value = 4161

# Row 34 data
$ws.Range("B34").Value = 443
$ws.Range("C34").Value = 188
$ws.Range("D34").Formula = "=B34-C34"
$ws.Range("E34").Value = 6
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 146.38
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0

# Final selection on the active cell
$ws.Range("N34").Select()
